# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" (column D) for the most recently
# handed-off file (row 5, fc98afc9-...) on both the zh-cn and de-de
# language report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-22 13:27:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-22 13:27:38"
